$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.839.61'

# Row 3
$ws.Range('D3').Value = '2.918.31'
$ws.Range('E3').Value = '  -3.89%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.39'
$ws.Range('E5').Value = '  -1.72%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.33'
$ws.Range('E6').Value = '  -6.01%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('E8').Value = '  -2.90%  '

# Row 9
$ws.Range('D9').Value = '2.916.89'
$ws.Range('E9').Value = '  -3.74%  '

# Row 10
$ws.Range('E10').Value = '  +5.92%  '

# Row 11
$ws.Range('E11').Value = '  -4.28%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.448'
$ws.Range('E12').Value = '  -3.82%  '

# Row 13
$ws.Range('E13').Value = '  -3.77%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.53'
$ws.Range('E14').Value = '  -5.50%  '

# Row 15
$ws.Range('E15').Value = '  +0.56%  '

# Row 16
$ws.Range('D16').Value = '3.400.70'
$ws.Range('E16').Value = '  -3.90%  '

# Row 17
$ws.Range('D17').Value = '60.764.42'
$ws.Range('E17').Value = '  -3.40%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.75'
$ws.Range('E18').Value = '  -4.65%  '

# Row 19
$ws.Range('D19').Value = '2.915.70'
$ws.Range('E19').Value = '  -3.98%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '431.35'
$ws.Range('E20').Value = '  -4.65%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.66'
$ws.Range('E21').Value = '  -4.20%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.683'
$ws.Range('E22').Value = '  -1.74%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.14'
$ws.Range('E23').Value = '  -4.60%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.54'
$ws.Range('E24').Value = '  -3.19%  '

# Row 25
$ws.Range('E25').Value = '  -1.71%  '

# Row 26
$ws.Range('E26').Value = '  -4.87%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.88'
$ws.Range('E27').Value = '  -3.08%  '

# Row 28
$ws.Range('E28').Value = '  +0.03%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.06%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.20'
$ws.Range('E30').Value = '  -4.87%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.62'
$ws.Range('E31').Value = '  -3.08%  '

# Row 32
$ws.Range('E32').Value = '  -3.14%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.57'
$ws.Range('E33').Value = '  -3.82%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.107'
$ws.Range('E34').Value = '  -4.04%  '

# Row 35
$ws.Range('D35').Value = '0.0₃0873'
$ws.Range('E35').Value = '  -0.14%  '

# Row 36
$ws.Range('E36').Value = '  -2.92%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.65'
$ws.Range('E37').Value = '  -4.63%  '

# Row 38
$ws.Range('E38').Value = '  -5.02%  '

# Row 39
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.127'
$ws.Range('E39').Value = '  +0.18%  '

# Row 40
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '49.84'
$ws.Range('E40').Value = '  -1.68%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.01'
$ws.Range('E41').Value = '  -4.42%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.70'
$ws.Range('E42').Value = '  -4.25%  '

# Row 43
$ws.Range('E43').Value = '  -5.20%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.89'
$ws.Range('E44').Value = '  -1.70%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '378.72'
$ws.Range('E45').Value = '  -4.41%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0348'
$ws.Range('E46').Value = '  -3.36%  '

# Row 47
$ws.Range('D47').Value = '2.678.29'
$ws.Range('E47').Value = '  -2.30%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.56'
$ws.Range('E48').Value = '  +0.09%  '

# Row 50
$ws.Range('E50').Value = '  +0.51%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.107'
$ws.Range('E51').Value = '  -1.62%  '
